$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1560865643779764
$ws.Range("C2").Value = 0.3515450347245845
$ws.Range("D2").Value = 0.2150214299408537
$ws.Range("E2").Value = 0.4637040326985023
$ws.Range("F2").Value = 0.446457955381491
$ws.Range("G2").Value = 23

$ws.Range("B3").Value = 0.03563623385673412
$ws.Range("C3").Value = 0.3763841266634617
$ws.Range("D3").Value = 0.258718923892501
$ws.Range("E3").Value = 0.5086442016699895
$ws.Range("F3").Value = 0.5193346301731255
$ws.Range("G3").Value = 22

$ws.Range("B4").Value = 0.1674435058222006
$ws.Range("C4").Value = 0.325651010243966
$ws.Range("D4").Value = 0.1893439557694782
$ws.Range("E4").Value = 0.4351367092874125
$ws.Range("F4").Value = 0.4115482469089393
$ws.Range("G4").Value = 21

$ws.Range("B5").Value = 0.09420853156004859
$ws.Range("C5").Value = 0.3793173975961993
$ws.Range("D5").Value = 0.1861955740678007
$ws.Range("E5").Value = 0.4315038517415583
$ws.Range("F5").Value = 0.4320335350642416
$ws.Range("G5").Value = 20

$ws.Range("B6").Value = 0.1376554976156691
$ws.Range("C6").Value = 0.3070384954770793
$ws.Range("D6").Value = 0.1917667010867467
$ws.Range("E6").Value = 0.4379117503410324
$ws.Range("F6").Value = 0.4271049595302239
$ws.Range("G6").Value = 19

$ws.Range("B7").Value = 0.08704673532814769
$ws.Range("C7").Value = 0.3331907123056521
$ws.Range("D7").Value = 0.1661864209314917
$ws.Range("E7").Value = 0.4076596876458251
$ws.Range("F7").Value = 0.4098039102391214
$ws.Range("G7").Value = 18

$ws.Range("B8").Value = 0.1098414975704748
$ws.Range("C8").Value = 0.3502453509676531
$ws.Range("D8").Value = 0.1720447973870636
$ws.Range("E8").Value = 0.4147828315963229
$ws.Range("F8").Value = 0.4122843320736888
$ws.Range("G8").Value = 17
